# The "pandas_experiment" sheet originally held a small A1:D4 table that was
# produced by seeking the wrong header row when reading the dataframe back
# in. Re-anchor the table at B3:F7 (header row 3, data rows 5-7, blank
# separator row 4, blank separator column D) with the corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("pandas_experiment")

# Clear out the old A1:D4 table entirely.
$ws.Range("A1:D4").ClearContents()

# Header row (row 3).
$ws.Range("B3").Value = "hello"

# Second little table's header cells (columns E/F) get filled first so the
# shared-string table lands in the same order Excel produced.
$ws.Range("E3").Value = "bye"
$ws.Range("F3").Value = "bye"

# Data rows for the first table (columns B/C).
$ws.Range("B5").Value = 1
$ws.Range("C5").Value = "a"

# Data rows for the second table (columns E/F), written in the order that
# matches the authored shared-string sequence.
$ws.Range("E5").Value = "fdsa"
$ws.Range("F7").Value = "fsda"
$ws.Range("F6").Value = "fds"
$ws.Range("F5").Value = "jona"
$ws.Range("E6").Value = "nic"
$ws.Range("E7").Value = "kri"

# Header for column C, filled in after the E/F data so it becomes the last
# new shared string.
$ws.Range("C3").Value = "sappy"

# Remaining first-table data cells.
$ws.Range("B6").Value = 2
$ws.Range("C6").Value = "b"
$ws.Range("B7").Value = 3
$ws.Range("C7").Value = "c"

# Match the saved selection/active cell.
$ws.Range("C4").Select()
